# Scheduled market-price refresh for the Leve profit tracker.
# Re-pulls current Market Board averages per item and recomputes the
# NQ/HQ totals + profit columns (H:N) for the leves whose prices moved.
$wb = $excel.ActiveWorkbook

# Updated figures: sheet name, row, column index, new value
$priceUpdates = @(
    ,@("ALC", 97, 8, 2825)
    ,@("ALC", 97, 10, 2825)
    ,@("ALC", 97, 12, 8475)
    ,@("ALC", 97, 14, -9467)
    ,@("ALC", 113, 8, 33337394)
    ,@("ALC", 113, 9, 55558456)
    ,@("ALC", 113, 11, 55558456)
    ,@("ALC", 113, 13, -55555202)
    ,@("ALC", 131, 9, 715171.9)
    ,@("ALC", 131, 10, 4099.5)
    ,@("ALC", 131, 11, 2145515.7)
    ,@("ALC", 131, 12, 12298.5)
    ,@("ALC", 131, 13, -2140475.7)
    ,@("ALC", 131, 14, -22378.5)
    ,@("ALC", 132, 8, 4837.9033)
    ,@("ALC", 132, 9, 5414.1113)
    ,@("ALC", 132, 10, 948.5)
    ,@("ALC", 132, 11, 16242.3339)
    ,@("ALC", 132, 12, 2845.5)
    ,@("ALC", 132, 13, -13712.3339)
    ,@("ALC", 132, 14, -7905.5)
    ,@("ALC", 135, 8, 3172)
    ,@("ALC", 135, 9, 1954.3334)
    ,@("ALC", 135, 11, 17589.0006)
    ,@("ALC", 135, 13, -15054.0006)
    ,@("ALC", 136, 8, 84174.5)
    ,@("ALC", 136, 10, 83663)
    ,@("ALC", 136, 12, 83663)
    ,@("ALC", 136, 14, -93863)
    ,@("ALC", 137, 8, 1725908.2)
    ,@("ALC", 137, 9, 2382812.2)
    ,@("ALC", 137, 11, 7148436.600000001)
    ,@("ALC", 137, 13, -7145886.600000001)
    ,@("ALC", 141, 8, 2735.9285)
    ,@("ALC", 141, 9, 2484.8462)
    ,@("ALC", 141, 10, 6000)
    ,@("ALC", 141, 11, 7454.5386)
    ,@("ALC", 141, 12, 18000)
    ,@("ALC", 141, 13, -2274.5386)
    ,@("ALC", 141, 14, -28360)
    ,@("ARM", 74, 8, 1823.2727)
    ,@("ARM", 74, 9, 1200)
    ,@("ARM", 74, 11, 1200)
    ,@("ARM", 74, 13, -326)
    ,@("ARM", 77, 8, 1823.2727)
    ,@("ARM", 77, 9, 1200)
    ,@("ARM", 77, 11, 6000)
    ,@("ARM", 77, 13, -1632)
    ,@("ARM", 80, 8, 71552.5)
    ,@("ARM", 80, 10, 94055)
    ,@("ARM", 80, 12, 94055)
    ,@("ARM", 80, 14, -96051)
    ,@("ARM", 83, 8, 71552.5)
    ,@("ARM", 83, 10, 94055)
    ,@("ARM", 83, 12, 282165)
    ,@("ARM", 83, 14, -292149)
    ,@("ARM", 109, 8, 75000)
    ,@("ARM", 109, 10, 75000)
    ,@("ARM", 109, 12, 75000)
    ,@("ARM", 109, 14, -77774)
    ,@("ARM", 132, 8, 2558.8572)
    ,@("ARM", 132, 9, 1304)
    ,@("ARM", 132, 10, 3500)
    ,@("ARM", 132, 11, 3912)
    ,@("ARM", 132, 12, 10500)
    ,@("ARM", 132, 13, -1382)
    ,@("ARM", 132, 14, -15560)
    ,@("BSM", 8, 8, 550)
    ,@("BSM", 8, 10, 550)
    ,@("BSM", 8, 12, 550)
    ,@("BSM", 8, 14, -830)
    ,@("BSM", 134, 8, 2883.125)
    ,@("BSM", 134, 9, 1966.375)
    ,@("BSM", 134, 10, 3799.875)
    ,@("BSM", 134, 11, 5899.125)
    ,@("BSM", 134, 12, 11399.625)
    ,@("BSM", 134, 13, -3364.125)
    ,@("BSM", 134, 14, -16469.625)
    ,@("CRP", 22, 8, 766.6667)
    ,@("CRP", 22, 9, 650)
    ,@("CRP", 22, 11, 650)
    ,@("CRP", 22, 13, -300)
    ,@("CRP", 31, 8, 5956164)
    ,@("CRP", 31, 9, 3518.5)
    ,@("CRP", 31, 10, 13893025)
    ,@("CRP", 31, 11, 3518.5)
    ,@("CRP", 31, 12, 13893025)
    ,@("CRP", 31, 13, -3223.5)
    ,@("CRP", 31, 14, -13893615)
    ,@("CRP", 34, 8, 5956164)
    ,@("CRP", 34, 9, 3518.5)
    ,@("CRP", 34, 10, 13893025)
    ,@("CRP", 34, 11, 3518.5)
    ,@("CRP", 34, 12, 13893025)
    ,@("CRP", 34, 13, -3316.5)
    ,@("CRP", 34, 14, -13893429)
    ,@("CRP", 58, 8, 1697.5416)
    ,@("CRP", 58, 9, 1377.5)
    ,@("CRP", 58, 11, 1377.5)
    ,@("CRP", 58, 13, -1174.5)
    ,@("CRP", 132, 8, 4055.1714)
    ,@("CRP", 132, 9, 3487.6667)
    ,@("CRP", 132, 11, 10463.0001)
    ,@("CRP", 132, 13, -7933.000100000001)
    ,@("CRP", 136, 8, 1697.5416)
    ,@("CRP", 136, 9, 1377.5)
    ,@("CRP", 136, 11, 4132.5)
    ,@("CRP", 136, 13, -1582.5)
    ,@("CRP", 140, 8, 69999)
    ,@("CRP", 140, 10, 69999)
    ,@("CRP", 140, 12, 69999)
    ,@("CRP", 140, 14, -80359)
    ,@("CUL", 2, 8, 138.425)
    ,@("CUL", 2, 9, 161.7)
    ,@("CUL", 2, 10, 115.15)
    ,@("CUL", 2, 11, 970.1999999999999)
    ,@("CUL", 2, 12, 690.9000000000001)
    ,@("CUL", 2, 13, -857.1999999999999)
    ,@("CUL", 2, 14, -916.9000000000001)
    ,@("CUL", 4, 8, 1680820.5)
    ,@("CUL", 4, 9, 1332642.5)
    ,@("CUL", 4, 11, 3997927.5)
    ,@("CUL", 4, 13, -3997815.5)
    ,@("CUL", 11, 8, 630)
    ,@("CUL", 11, 9, 606.6667)
    ,@("CUL", 11, 10, 700)
    ,@("CUL", 11, 11, 1820.0001)
    ,@("CUL", 11, 12, 2100)
    ,@("CUL", 11, 13, -1680.0001)
    ,@("CUL", 11, 14, -2380)
    ,@("CUL", 107, 8, 475.07693)
    ,@("CUL", 107, 9, 382)
    ,@("CUL", 107, 10, 516.44446)
    ,@("CUL", 107, 11, 1146)
    ,@("CUL", 107, 12, 1549.33338)
    ,@("CUL", 107, 13, 774)
    ,@("CUL", 107, 14, -5389.33338)
    ,@("CUL", 132, 8, 1251.6666)
    ,@("CUL", 132, 10, 1200)
    ,@("CUL", 132, 12, 10800)
    ,@("CUL", 132, 14, -15860)
    ,@("GSM", 102, 8, 4466.8237)
    ,@("GSM", 102, 9, 1735)
    ,@("GSM", 102, 11, 1735)
    ,@("GSM", 102, 13, -113)
    ,@("GSM", 132, 8, 1937.9131)
    ,@("GSM", 132, 9, 1452)
    ,@("GSM", 132, 10, 2383.3333)
    ,@("GSM", 132, 11, 4356)
    ,@("GSM", 132, 12, 7149.999899999999)
    ,@("GSM", 132, 13, -1826)
    ,@("GSM", 132, 14, -12209.9999)
    ,@("LTW", 46, 8, 2473.3)
    ,@("LTW", 46, 9, 1942.7142)
    ,@("LTW", 46, 10, 3711.3333)
    ,@("LTW", 46, 11, 1942.7142)
    ,@("LTW", 46, 12, 3711.3333)
    ,@("LTW", 46, 13, -1754.7142)
    ,@("LTW", 46, 14, -4087.3333)
    ,@("LTW", 56, 8, 12220.223)
    ,@("LTW", 56, 10, 26249)
    ,@("LTW", 56, 12, 26249)
    ,@("LTW", 56, 14, -27631)
    ,@("LTW", 64, 8, 30150)
    ,@("LTW", 64, 10, 30150)
    ,@("LTW", 64, 12, 30150)
    ,@("LTW", 64, 14, -30600)
    ,@("LTW", 67, 8, 30150)
    ,@("LTW", 67, 10, 30150)
    ,@("LTW", 67, 12, 30150)
    ,@("LTW", 67, 14, -31710)
    ,@("LTW", 132, 8, 3387.5386)
    ,@("LTW", 132, 9, 3221.5557)
    ,@("LTW", 132, 10, 3761)
    ,@("LTW", 132, 11, 9664.667099999999)
    ,@("LTW", 132, 12, 11283)
    ,@("LTW", 132, 13, -7134.667099999999)
    ,@("LTW", 132, 14, -16343)
    ,@("LTW", 136, 8, 4987.375)
    ,@("LTW", 136, 9, 5089.316)
    ,@("LTW", 136, 11, 15267.948)
    ,@("LTW", 136, 13, -12717.948)
    ,@("LTW", 140, 8, 118906.91)
    ,@("LTW", 140, 10, 118906.91)
    ,@("LTW", 140, 12, 118906.91)
    ,@("LTW", 140, 14, -129266.91)
    ,@("WVR", 34, 8, 10000)
    ,@("WVR", 34, 9, 0)
    ,@("WVR", 34, 10, 10000)
    ,@("WVR", 34, 11, 0)
    ,@("WVR", 34, 12, 10000)
    ,@("WVR", 34, 14, -10406)
    ,@("WVR", 70, 8, 0)
    ,@("WVR", 70, 10, 0)
    ,@("WVR", 70, 12, 0)
    ,@("WVR", 73, 8, 0)
    ,@("WVR", 73, 10, 0)
    ,@("WVR", 73, 12, 0)
    ,@("WVR", 76, 8, 28500)
    ,@("WVR", 76, 10, 28500)
    ,@("WVR", 76, 12, 28500)
    ,@("WVR", 76, 14, -29130)
    ,@("WVR", 79, 8, 28500)
    ,@("WVR", 79, 10, 28500)
    ,@("WVR", 79, 12, 28500)
    ,@("WVR", 79, 14, -30684)
    ,@("WVR", 132, 8, 2797)
    ,@("WVR", 132, 9, 2770.0386)
    ,@("WVR", 132, 10, 2937.2)
    ,@("WVR", 132, 11, 8310.1158)
    ,@("WVR", 132, 12, 8811.599999999999)
    ,@("WVR", 132, 13, -5780.1158)
    ,@("WVR", 132, 14, -13871.6)
    ,@("WVR", 136, 8, 5700.857)
    ,@("WVR", 136, 9, 5700.857)
    ,@("WVR", 136, 10, 0)
    ,@("WVR", 136, 11, 17102.571)
    ,@("WVR", 136, 12, 0)
    ,@("WVR", 136, 13, -14552.571)
)

# Columns that no longer apply for these leves (item now only sells one way)
# and must be cleared out rather than recomputed
$staleColumns = @(
    ,@("WVR", 34, 13)
    ,@("WVR", 70, 14)
    ,@("WVR", 136, 14)
    ,@("WVR", 73, 14)
)

foreach ($u in $priceUpdates) {
    $sheetName = $u[0]
    $row = $u[1]
    $col = $u[2]
    $val = $u[3]
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item($row, $col).Value = $val
}

foreach ($c in $staleColumns) {
    $sheetName = $c[0]
    $row = $c[1]
    $col = $c[2]
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item($row, $col).ClearContents()
}

Write-Output "Applied $($priceUpdates.Count) value updates and $($staleColumns.Count) clears."
